# Apply updated cryptocurrency price/volume figures (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.200.19"
$ws.Range("E2").Value = "  -1.40%  "
$ws.Range("D3").Value = "3.396.26"
$ws.Range("E3").Value = "  -1.48%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'569.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.21%  "
$ws.Range("D6").Value = "'155.72"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.24%  "
$ws.Range("D7").Value = "'0.628"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +8.62%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "3.399.34"
$ws.Range("D10").Value = "'7.12"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.91%  "
$ws.Range("E11").Value = "  -2.25%  "
$ws.Range("D13").Value = "3.984.92"
$ws.Range("E13").Value = "  -1.46%  "
$ws.Range("E14").Value = "  -0.30%  "
$ws.Range("E15").Value = "  -3.49%  "
$ws.Range("D16").Value = "'27.33"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.55%  "
$ws.Range("D17").Value = "64.234.74"
$ws.Range("E17").Value = "  -1.35%  "
$ws.Range("D18").Value = "3.403.40"
$ws.Range("E18").Value = "  -1.33%  "
$ws.Range("D19").Value = "'6.27"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.08%  "
$ws.Range("D20").Value = "'13.76"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.30%  "
$ws.Range("D21").Value = "'376.53"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.79%  "
$ws.Range("D22").Value = "'7.95"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.85%  "
$ws.Range("E24").Value = "  -0.36%  "
$ws.Range("D25").Value = "'71.63"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.11%  "
$ws.Range("D26").Value = "'0.0000117"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.44%  "
$ws.Range("D27").Value = "'10.32"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +7.02%  "
$ws.Range("E28").Value = "  -2.23%  "
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("E30").Value = "  +1.92%  "
$ws.Range("E31").Value = "  -2.49%  "
$ws.Range("E32").Value = "  -1.92%  "
$ws.Range("D33").Value = "'23.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.50%  "
$ws.Range("D34").Value = "'7.07"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.24%  "
$ws.Range("D35").Value = "'1.58"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.95%  "
$ws.Range("D36").Value = "'160.33"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.92%  "
$ws.Range("D37").Value = "'1.88"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.72%  "
$ws.Range("D38").Value = "'0.0755"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.43%  "
$ws.Range("D39").Value = "2.864.05"
$ws.Range("E39").Value = "  -5.79%  "
$ws.Range("D40").Value = "'6.74"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.13%  "
$ws.Range("D41").Value = "'26.22"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.72%  "
$ws.Range("D42").Value = "'4.56"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.55%  "
$ws.Range("D43").Value = "'42.60"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.39%  "
$ws.Range("E44").Value = "  -1.15%  "
$ws.Range("D45").Value = "'0.766"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.52%  "
$ws.Range("D46").Value = "'25.69"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.43%  "
$ws.Range("D47").Value = "'321.31"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.36%  "
$ws.Range("E48").Value = "  +3.33%  "
$ws.Range("E49").Value = "  -2.49%  "
$ws.Range("D50").Value = "'2.17"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.30%  "
$ws.Range("D51").Value = "'6.50"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.53%  "
